# New features of alphabetical sort, input focus, and tag suggestions.
# This populates the demo/sample data used to exercise those features:
#   - row 2 gets a run of single-letter "tag" values across D2:J2
#   - rows 3-5 get additional "tag" / suggestion sample values in B:E
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (D2:J2) - alphabetical-sort / tag-suggestion sample tags
$ws.Range("D2").Value = "a"
$ws.Range("E2").Value = "c"
$ws.Range("F2").Value = "b"
$ws.Range("G2").Value = "x"
$ws.Range("H2").Value = "e"
$ws.Range("I2").Value = "f"
$ws.Range("J2").Value = "g"

# Row 3 (B3:C3)
$ws.Range("B3").Value = "new"
$ws.Range("C3").Value = "hi"

# Row 4 (B4:I4)
$ws.Range("B4").Value = "hello"
$ws.Range("C4").Value = "hi"
$ws.Range("D4").Value = "new"
$ws.Range("E4").Value = "hey"
$ws.Range("F4").Value = "help"
$ws.Range("G4").Value = "heyow"
$ws.Range("H4").Value = "helpo"
$ws.Range("I4").Value = "helpo"

# Row 5 (D5:E5)
$ws.Range("D5").Value = "new tag"
$ws.Range("E5").Value = "banana"
